$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A values (rows 2-11) from 1 to 9040000000
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = 9040000000
}

# Update the active selection to D1
$ws.Range("D1").Select()
